$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix up the note text for the prior entry and add two new shared-string rows
$ws.Range("A36").Value = "2012.6.28"
$ws.Range("B36").Value = "修改一处属性球bug，判断球回收为hitpoint为-1"
$ws.Range("D36").Value = 2

$ws.Range("A36").Style = $ws.Range("A35").Style
$ws.Range("B36").Style = $ws.Range("B35").Style
$ws.Rows.Item(36).RowHeight = $ws.Rows.Item(32).RowHeight

$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("A22").Select()
